$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is missing a "recovery" email in column D. Add it back as a
# mailto hyperlink, matching the same pattern used by the other
# recovery-email cells in column D (e.g. D4, D8, D9, D10, D13).
[void]$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:budoyoni@gmail.com", "", "", "budoyoni@gmail.com")

# Match the plain (non-auto-hyperlink) formatting used by the sibling
# recovery-email cells: Arial 10, blue text, no underline.
$srcFont = $ws.Range("D4").Font
$dstFont = $ws.Range("D3").Font
$dstFont.Name = $srcFont.Name()
$dstFont.Size = $srcFont.Size()
$dstFont.Color = $srcFont.Color()
$dstFont.Underline = $srcFont.Underline()

# Leave the cursor where the author left it after the edit.
[void]$ws.Range("D4").Select()
